# تعديل تلقائي في شيت Card22
#
# The source export re-ran its pandas->xlsx dump: the "Card3" lookup table
# had its literal "nan" placeholder strings stripped back down to blank
# cells, while the "Card22" lookup table had its blanks re-filled with the
# "nan" placeholder text and gained the trailing (empty) row that every
# other card sheet in this workbook already carries.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Card3": drop the literal "nan" text, leaving the cells blank.
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Card3")

$ws3.Range("D2:N2").ClearContents()

$ws3.Range("D3").ClearContents()
$ws3.Range("F3:K3").ClearContents()
$ws3.Range("M3:N3").ClearContents()

$ws3.Range("D4:N4").ClearContents()

$ws3.Range("D5").ClearContents()
$ws3.Range("H5").ClearContents()
$ws3.Range("J5:K5").ClearContents()
$ws3.Range("M5:N5").ClearContents()

$ws3.Range("E6:G6").ClearContents()
$ws3.Range("I6:K6").ClearContents()
$ws3.Range("M6:N6").ClearContents()

$ws3.Range("E7").ClearContents()
$ws3.Range("G7:J7").ClearContents()
$ws3.Range("M7:N7").ClearContents()

$ws3.Range("F8:K8").ClearContents()

$ws3.Range("E9").ClearContents()
$ws3.Range("H9:K9").ClearContents()
$ws3.Range("M9:N9").ClearContents()

$ws3.Range("D10:N10").ClearContents()
$ws3.Range("D11:N11").ClearContents()
$ws3.Range("D12:N12").ClearContents()
$ws3.Range("D13:N13").ClearContents()

# ------------------------------------------------------------------
# Sheet "Card22": restore the "nan" placeholder text in the blank data
# cells, and append the trailing blank row 13.
# ------------------------------------------------------------------
$ws22 = $wb.Worksheets.Item("Card22")

# --- Row 2 ---
$ws22.Range("E2:K2").Value = "nan"
$ws22.Range("N2").Value = "nan"

# --- Row 3 ---
$ws22.Range("D3:O3").Value = "nan"

# --- Row 4 ---
$ws22.Range("D4:O4").Value = "nan"

# --- Row 5 ---
$ws22.Range("D5").Value = "nan"
$ws22.Range("H5").Value = "nan"
$ws22.Range("J5:K5").Value = "nan"
$ws22.Range("M5:N5").Value = "nan"

# --- Row 6 ---
$ws22.Range("E6:G6").Value = "nan"
$ws22.Range("I6:K6").Value = "nan"
$ws22.Range("M6:O6").Value = "nan"

# --- Row 7 ---
$ws22.Range("E7").Value = "nan"
$ws22.Range("H7:K7").Value = "nan"
$ws22.Range("M7:O7").Value = "nan"

# --- Row 8 ---
$ws22.Range("D8:O8").Value = "nan"

# --- Row 9 ---
$ws22.Range("D9:O9").Value = "nan"

# --- Row 10 ---
$ws22.Range("D10:O10").Value = "nan"

# --- Row 11 ---
$ws22.Range("D11:O11").Value = "nan"

# --- Row 12 ---
$ws22.Range("D12:O12").Value = "nan"

# --- New trailing blank row 13 (extends used range without adding values) ---
$ws22.Range("A13:O13").Style = $ws22.Range("A12").Style
